$d = $word.ActiveDocument

# 1. Merge the two runs "(对于图片展示" + [bookmark] + "的控制)" into a single
#    run "(对于图片展示的控制)". This also removes the old _GoBack bookmark
#    that used to sit between them (it is inside the replaced range).
$null = $d.Content.Find.Execute("(对于图片展示的控制)", $true, $false, $false, $false, $false, $true, 1, $false, "(对于图片展示的控制)", 2)

# 2. Re-create the _GoBack bookmark right after "参考文献：" -- collapsed,
#    zero-width, immediately before that paragraph's mark.
$findRng = $d.Content
$null = $findRng.Find.Execute("参考文献：", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$para = $findRng.Paragraphs(1)
$r = $para.Range
$r.Collapse(0)
$null = $r.MoveEnd(1, -1)

# Inserting a marker character lets us create a non-collapsed range that the
# COM host anchors correctly, then we delete the marker leaving a correctly
# placed, zero-width bookmark.
$r.InsertAfter("Z")
$d.Bookmarks.Add("_GoBack", $r)
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Delete()
